$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the 6 additional entrance rows (rows 8-13 are brand new).
$ws.Rows("8:13").Insert()

# ---- Column C: the entrance names (plain text) for every data row ----
$names = @{
    1  = "Wing D North Parking Garage Entrance"
    2  = "South Limestone Entrance A"
    3  = "South Limestone Entrance B"
    4  = "Rose Street Entrance"
    5  = "University Health Services Entrance"
    6  = "1st Floor Main Entrance"
    7  = "Charles T. Wethington Building Entrance"
    8  = "2nd Floor Main Entrance"
    9  = "University Health Services Bridge Entrance"
    10 = "A.B. Chandler Bridge Entrance"
    11 = "3rd Floor Parking Garage Entrance"
    12 = "4th Floor Parking Garage Entrance"
    13 = "5th Floor Parking Garage Entrance"
}

foreach ($r in $names.Keys) {
    $ws.Cells.Item($r, 3).Value = $names[$r]
}

# ---- Column B: filename formula, =CONCAT(SUBSTITUTE(C#," ","_"),".png") ----
# Rows 4 and 5 ended up pointing at the wrong source cell (C2 / C4 respectively)
# after the sheet was edited - reproduce that exactly.
$bRefs = @{
    1  = "C1"
    2  = "C2"
    3  = "C3"
    4  = "C2"
    5  = "C4"
    6  = "C6"
    7  = "C7"
    8  = "C8"
    9  = "C9"
    10 = "C10"
    11 = "C11"
    12 = "C12"
    13 = "C13"
}

foreach ($r in $bRefs.Keys) {
    $ws.Cells.Item($r, 2).Formula = '=_xlfn.CONCAT(SUBSTITUTE(' + $bRefs[$r] + ',' + '" "' + ',"_"),".png")'
}

# ---- Column A: the start-url formula, always references same-row B ----
for ($r = 1; $r -le 13; $r++) {
    $ws.Cells.Item($r, 1).Formula = '=_xlfn.CONCAT("https://hospitalnavigation.netlify.app/start/?start=",SUBSTITUTE(B' + $r + ',".png",""))'
}

# ---- Stray leftover cell far below the table ----
$ws.Cells.Item(24, 1).Value = "s"

# Restore the on-screen selection to the populated block.
$ws.Range("A1:C13").Select()
